$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Garsoniera, 32mp, cartier Manastur, zon Big"
$ws.Range("B1").Value = "230 EUR / lună"
$ws.Range("C1").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/manastur/apartament-de-inchiriat-2-camere-X5NQ1022U"

$ws.Range("A2").Value = "Apartament 2 camere, recent renovat, 37mp, Manastur , Pet friendly"
$ws.Range("B2").Value = "240 EUR / lună"
$ws.Range("C2").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/manastur/apartament-de-inchiriat-2-camere-X3SH10DA2"

$ws.Range("A3").Value = "Apartament 2 camere , modest , mobilat-utilat!!!"
$ws.Range("B3").Value = "250 EUR / lună"
$ws.Range("C3").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/manastur/apartament-de-inchiriat-2-camere-XAON100FG"

$ws.Range("A4").Value = "Apartament 2 camere, decomandat, 45 mp, pet friendly, zona strazii..."
$ws.Range("B4").Value = "270 EUR / lună"
$ws.Range("C4").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/manastur/apartament-de-inchiriat-2-camere-X3SH10DDL"

$ws.Range("A5").Value = "Apartament 2 camere Manastur, str. Tasnad, zona Cora, partial mobilat, cu GARAJ"
$ws.Range("B5").Value = "270 EUR / lună"
$ws.Range("C5").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/manastur/apartament-de-inchiriat-2-camere-XARU0003Q"

# Preserve the pre-existing empty row 11 (leftover selection artifact) so it
# is not dropped from sheetData when the used range is recalculated.
$ws.Cells.Item(11, 1).EntireRow.OutlineLevel = 0
